$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell K1: copy the existing header formatting (style) from J1,
# then set its text to "Request Frequency".
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("K1").Value = "Request Frequency"

# New data column values
$ws.Range("K2").Value = 0.8
$ws.Range("K3").Value = 0.45
